$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by 2 days
# (from 2023-09-06 / serial 45175 to 2023-09-08 / serial 45177) for every
# data row (rows 2 through 484).
$rng = $ws.Range("C2:C484")
$rng.Value = 45177
